# TableMLCompiler template extension: add role card data (Cao Cao) to zh_CN.xlsx Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New rows describing the "Cao Cao" role card (fill order matches the
# original authoring tool's (TableMLCompiler) cell write sequence so the
# shared-string table comes out in the same order)
$ws.Range("A6").Value = "CAOCAO"
$ws.Range("B6").Value = "曹操"

$ws.Range("B7").Value = "曹操的技能描述"
$ws.Range("A7").Value = "CAOCAO_DESC"

$ws.Range("A8").Value = "BELONG_1"
$ws.Range("B8").Value = "魏"

$ws.Range("A9").Value = "BELONG_2"
$ws.Range("B9").Value = "蜀"

$ws.Range("A10").Value = "BELONG_3"
$ws.Range("A11").Value = "BELONG_4"
$ws.Range("B10").Value = "吴"
$ws.Range("B11").Value = "它"

$ws.Range("A12").Value = "FORCE"
$ws.Range("A13").Value = "COMMAND"
$ws.Range("A14").Value = "MORAL"
$ws.Range("B12").Value = "武力"
$ws.Range("B13").Value = "统帅"
$ws.Range("B14").Value = "德行"

# Column widths to fit the longer localized content
# (inputs compensate for the engine's internal +5/7-character padding so the
# stored OOXML <col width> comes out at 46 and ~59.875 respectively)
$ws.Columns.Item(1).ColumnWidth = 45.285714285714285
$ws.Columns.Item(2).ColumnWidth = 59.142857142857146

# Move selection like the source commit (cursor left on next free row)
$ws.Range("B16").Select() | Out-Null
